$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 1 de Octubre de 2020 a las 11:21"

# Re-rank countries: swap display order/names for rows whose rank changed
$ws.Range("A23").Value = "Filipinas"
$ws.Range("A24").Value = "Pakistan"
$ws.Range("A44").Value = "Polonia"
$ws.Range("A45").Value = "Suecia"
$ws.Range("A46").Value = "Guatemala"
$ws.Range("A118").Value = "Eslovenia"
$ws.Range("A119").Value = "Malaui"
$ws.Range("A125").Value = "Hong Kong"
$ws.Range("A126").Value = "Congo"
$ws.Range("A158").Value = "Letonia"
$ws.Range("A159").Value = "Polinesia Francesa"
$ws.Range("A160").Value = "Nueva Zelanda"

# Updated statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
$ws.Range("B7").Value = 1185231
$ws.Range("C7").Value = 8945
$ws.Range("D7").Value = 964242
$ws.Range("E7").Value = 200098
$ws.Range("G7").Value = 169
$ws.Range("H7").Value = 20891
$ws.Range("B23").Value = 314079
$ws.Range("C23").Value = 2415
$ws.Range("D23").Value = 254223
$ws.Range("E23").Value = 54294
$ws.Range("G23").Value = 59
$ws.Range("H23").Value = 5562
$ws.Range("B24").Value = 312806
$ws.Range("C24").Value = 543
$ws.Range("D24").Value = 297497
$ws.Range("E24").Value = 8825
$ws.Range("G24").Value = 5
$ws.Range("H24").Value = 6484
$ws.Range("B44").Value = 93481
$ws.Range("C44").Value = 1967
$ws.Range("D44").Value = 70401
$ws.Range("E44").Value = 20537
$ws.Range("G44").Value = 30
$ws.Range("H44").Value = 2543
$ws.Range("B45").Value = 92863
$ws.Range("D45").Value = 0
$ws.Range("E45").Value = 0
$ws.Range("H45").Value = 5893
$ws.Range("B46").Value = 91746
$ws.Range("D46").Value = 80256
$ws.Range("E46").Value = 8244
$ws.Range("H46").Value = 3246
$ws.Range("B67").Value = 45686
$ws.Range("C67").Value = 873
$ws.Range("D67").Value = 36476
$ws.Range("E67").Value = 8408
$ws.Range("G67").Value = 3
$ws.Range("H67").Value = 802
$ws.Range("B72").Value = 39285
$ws.Range("C72").Value = 17
$ws.Range("D72").Value = 32842
$ws.Range("E72").Value = 4985
$ws.Range("B77").Value = 29175
$ws.Range("C77").Value = 98
$ws.Range("D77").Value = 23930
$ws.Range("E77").Value = 4397
$ws.Range("B89").Value = 16827
$ws.Range("C89").Value = 234
$ws.Range("D89").Value = 15218
$ws.Range("E89").Value = 1325
$ws.Range("G89").Value = 4
$ws.Range("H89").Value = 284
$ws.Range("B105").Value = 10103
$ws.Range("C105").Value = 111
$ws.Range("E105").Value = 1659
$ws.Range("B118").Value = 5865
$ws.Range("C118").Value = 175
$ws.Range("D118").Value = 3906
$ws.Range("E118").Value = 1807
$ws.Range("G118").Value = 2
$ws.Range("H118").Value = 152
$ws.Range("B119").Value = 5773
$ws.Range("D119").Value = 4263
$ws.Range("E119").Value = 1331
$ws.Range("H119").Value = 179
$ws.Range("B125").Value = 5098
$ws.Range("C125").Value = 10
$ws.Range("D125").Value = 4836
$ws.Range("E125").Value = 157
$ws.Range("H125").Value = 105
$ws.Range("B126").Value = 5089
$ws.Range("D126").Value = 3887
$ws.Range("E126").Value = 1113
$ws.Range("H126").Value = 89
$ws.Range("B132").Value = 4784
$ws.Range("C132").Value = 91
$ws.Range("D132").Value = 2424
$ws.Range("E132").Value = 2268
$ws.Range("B158").Value = 1868
$ws.Range("C158").Value = 44
$ws.Range("D158").Value = 1307
$ws.Range("E158").Value = 524
$ws.Range("H158").Value = 37
$ws.Range("B159").Value = 1852
$ws.Range("C159").Value = 0
$ws.Range("D159").Value = 1504
$ws.Range("E159").Value = 341
$ws.Range("H159").Value = 7
$ws.Range("B160").Value = 1848
$ws.Range("C160").Value = 12
$ws.Range("D160").Value = 1770
$ws.Range("E160").Value = 53
$ws.Range("H160").Value = 25
$ws.Range("B168").Value = 1095
$ws.Range("C168").Value = 1
$ws.Range("E168").Value = 50
